# Auto-generated Excel COM-interop script to apply the Ramuh_Profits.xlsx diff
# Updates columns H-N (currentAveragePrice.. LeveProfitHQ) for specific rows across the 8 job sheets
$wb = $excel.ActiveWorkbook

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2339.36
$ws.Range("I113").Value = 1564.6666
$ws.Range("K113").Value = 1564.6666
$ws.Range("M113").Value = 1689.3334

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2651.2307
$ws.Range("I116").Value = 2562.2
$ws.Range("J116").Value = 2706.875
$ws.Range("K116").Value = 2562.2
$ws.Range("L116").Value = 2706.875
$ws.Range("M116").Value = 879.8000000000002
$ws.Range("N116").Value = -9590.875

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2846.9119
$ws.Range("I132").Value = 2512.4092
$ws.Range("J132").Value = 3460.1667
$ws.Range("K132").Value = 7537.2276
$ws.Range("L132").Value = 10380.5001
$ws.Range("M132").Value = -5007.2276
$ws.Range("N132").Value = -15440.5001

# Sheet ALC, row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 345.46667
$ws.Range("I135").Value = 331.9535
$ws.Range("J135").Value = 636
$ws.Range("K135").Value = 2987.5815
$ws.Range("L135").Value = 5724
$ws.Range("M135").Value = -452.5815000000002
$ws.Range("N135").Value = -10794

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 125338.875
$ws.Range("I2").Value = 166968.5
$ws.Range("J2").Value = 450
$ws.Range("K2").Value = 166968.5
$ws.Range("L2").Value = 450
$ws.Range("M2").Value = -166855.5
$ws.Range("N2").Value = -676

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2060.8667
$ws.Range("I61").Value = 1900.2858
$ws.Range("J61").Value = 2201.375
$ws.Range("K61").Value = 1900.2858
$ws.Range("L61").Value = 2201.375
$ws.Range("M61").Value = -1688.2858
$ws.Range("N61").Value = -2625.375

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 125338.875
$ws.Range("I116").Value = 166968.5
$ws.Range("J116").Value = 450
$ws.Range("K116").Value = 166968.5
$ws.Range("L116").Value = 450
$ws.Range("M116").Value = -164674.5
$ws.Range("N116").Value = -5038

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7036.4116
$ws.Range("I132").Value = 4761.5713
$ws.Range("J132").Value = 17652.334
$ws.Range("K132").Value = 14284.7139
$ws.Range("L132").Value = 52957.00199999999
$ws.Range("M132").Value = -11754.7139
$ws.Range("N132").Value = -58017.00199999999

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2060.8667
$ws.Range("I136").Value = 1900.2858
$ws.Range("J136").Value = 2201.375
$ws.Range("K136").Value = 5700.857400000001
$ws.Range("L136").Value = 6604.125
$ws.Range("M136").Value = -3150.857400000001
$ws.Range("N136").Value = -11704.125

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 125338.875
$ws.Range("I3").Value = 166968.5
$ws.Range("J3").Value = 450
$ws.Range("K3").Value = 166968.5
$ws.Range("L3").Value = 450
$ws.Range("M3").Value = -166854.5
$ws.Range("N3").Value = -678

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 83334680
$ws.Range("I99").Value = 125000740
$ws.Range("J99").Value = 2558.25
$ws.Range("K99").Value = 125000740
$ws.Range("L99").Value = 2558.25
$ws.Range("M99").Value = -124999242
$ws.Range("N99").Value = -5554.25

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 763.59155
$ws.Range("I134").Value = 677.31146
$ws.Range("J134").Value = 1289.9
$ws.Range("K134").Value = 2031.93438
$ws.Range("L134").Value = 3869.7
$ws.Range("M134").Value = 503.0656199999999
$ws.Range("N134").Value = -8939.700000000001

# Sheet CRP, row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3867.439
$ws.Range("I62").Value = 4756.087
$ws.Range("J62").Value = 2731.9443
$ws.Range("K62").Value = 4756.087
$ws.Range("L62").Value = 2731.9443
$ws.Range("M62").Value = -4132.087
$ws.Range("N62").Value = -3979.9443

# Sheet CRP, row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3867.439
$ws.Range("I65").Value = 4756.087
$ws.Range("J65").Value = 2731.9443
$ws.Range("K65").Value = 23780.435
$ws.Range("L65").Value = 13659.7215
$ws.Range("M65").Value = -20660.435
$ws.Range("N65").Value = -19899.7215

# Sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3760.5715
$ws.Range("I99").Value = 3904.8
$ws.Range("J99").Value = 3400
$ws.Range("K99").Value = 3904.8
$ws.Range("L99").Value = 3400
$ws.Range("M99").Value = -2406.8
$ws.Range("N99").Value = -6396

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 948.93335
$ws.Range("I107").Value = 929.1818
$ws.Range("J107").Value = 1003.25
$ws.Range("K107").Value = 929.1818
$ws.Range("L107").Value = 1003.25
$ws.Range("M107").Value = 990.8182
$ws.Range("N107").Value = -4843.25

# Sheet CRP, row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1302.875
$ws.Range("I122").Value = 1346.2858
$ws.Range("J122").Value = 999
$ws.Range("K122").Value = 4038.8574
$ws.Range("L122").Value = 2997
$ws.Range("M122").Value = -1588.8574
$ws.Range("N122").Value = -7897

# Sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3760.5715
$ws.Range("I126").Value = 3904.8
$ws.Range("J126").Value = 3400
$ws.Range("K126").Value = 11714.4
$ws.Range("L126").Value = 10200
$ws.Range("M126").Value = -9244.400000000001
$ws.Range("N126").Value = -15140

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 14087342
$ws.Range("I132").Value = 18521776
$ws.Range("J132").Value = 1495.6471
$ws.Range("K132").Value = 55565328
$ws.Range("L132").Value = 4486.9413
$ws.Range("M132").Value = -55562798
$ws.Range("N132").Value = -9546.941299999999

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4730.625
$ws.Range("I134").Value = 5113.4287
$ws.Range("J134").Value = 2051
$ws.Range("K134").Value = 15340.2861
$ws.Range("L134").Value = 6153
$ws.Range("M134").Value = -12805.2861
$ws.Range("N134").Value = -11223

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 263747.53
$ws.Range("I5").Value = 229.1875
$ws.Range("J5").Value = 455397.22
$ws.Range("K5").Value = 687.5625
$ws.Range("L5").Value = 1366191.66
$ws.Range("M5").Value = -575.5625
$ws.Range("N5").Value = -1366415.66

# Sheet CUL, row 36
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1234
$ws.Range("I36").Value = 1234
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 3702
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -3533
$ws.Range("N36").ClearContents()

# Sheet CUL, row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1751.6364
$ws.Range("I68").Value = 679.1667
$ws.Range("J68").Value = 3038.6
$ws.Range("K68").Value = 2037.5001
$ws.Range("L68").Value = 9115.799999999999
$ws.Range("M68").Value = -1226.5001
$ws.Range("N68").Value = -10737.8

# Sheet CUL, row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1751.6364
$ws.Range("I71").Value = 679.1667
$ws.Range("J71").Value = 3038.6
$ws.Range("K71").Value = 6112.5003
$ws.Range("L71").Value = 27347.4
$ws.Range("M71").Value = -2056.5003
$ws.Range("N71").Value = -35459.39999999999

# Sheet CUL, row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 452.5
$ws.Range("J107").Value = 468
$ws.Range("L107").Value = 1404
$ws.Range("N107").Value = -5244

# Sheet CUL, row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 695433.1
$ws.Range("I113").Value = 450.1579
$ws.Range("J113").Value = 1575744.9
$ws.Range("K113").Value = 1350.4737
$ws.Range("L113").Value = 4727234.699999999
$ws.Range("M113").Value = 819.5263
$ws.Range("N113").Value = -4731574.699999999

# Sheet CUL, row 116
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 4152.6665
$ws.Range("I116").Value = 5229
$ws.Range("K116").Value = 15687
$ws.Range("M116").Value = -12245

# Sheet CUL, row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1328.5714
$ws.Range("J117").Value = 1280
$ws.Range("L117").Value = 3840
$ws.Range("N117").Value = -10724

# Sheet CUL, row 119
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 250001250
$ws.Range("I119").Value = 250001250
$ws.Range("K119").Value = 750003750
$ws.Range("M119").Value = -749998912

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 263747.53
$ws.Range("I135").Value = 229.1875
$ws.Range("J135").Value = 455397.22
$ws.Range("K135").Value = 2062.6875
$ws.Range("L135").Value = 4098574.98
$ws.Range("M135").Value = 472.3125
$ws.Range("N135").Value = -4103644.98

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2610.2104
$ws.Range("J80").Value = 3080.3635
$ws.Range("L80").Value = 3080.3635
$ws.Range("N80").Value = -5076.363499999999

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2610.2104
$ws.Range("J83").Value = 3080.3635
$ws.Range("L83").Value = 15401.8175
$ws.Range("N83").Value = -25385.8175

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5948.815
$ws.Range("I132").Value = 7420.9
$ws.Range("J132").Value = 1742.8572
$ws.Range("K132").Value = 22262.7
$ws.Range("L132").Value = 5228.571599999999
$ws.Range("M132").Value = -19732.7
$ws.Range("N132").Value = -10288.5716

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 811.1875
$ws.Range("I93").Value = 654.36365
$ws.Range("J93").Value = 1156.2
$ws.Range("K93").Value = 654.36365
$ws.Range("L93").Value = 1156.2
$ws.Range("M93").Value = 593.63635
$ws.Range("N93").Value = -3652.2

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 50527.715
$ws.Range("I122").Value = 69108.13
$ws.Range("J122").Value = 4076.6667
$ws.Range("K122").Value = 207324.39
$ws.Range("L122").Value = 12230.0001
$ws.Range("M122").Value = -204874.39
$ws.Range("N122").Value = -17130.0001

# Sheet LTW, row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 38163
$ws.Range("J133").Value = 38163
$ws.Range("L133").Value = 38163
$ws.Range("N133").Value = -43223

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2739.2
$ws.Range("I62").Value = 2786
$ws.Range("J62").Value = 2630
$ws.Range("K62").Value = 2786
$ws.Range("L62").Value = 2630
$ws.Range("M62").Value = -2162
$ws.Range("N62").Value = -3878

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 2739.2
$ws.Range("I65").Value = 2786
$ws.Range("J65").Value = 2630
$ws.Range("K65").Value = 13930
$ws.Range("L65").Value = 13150
$ws.Range("M65").Value = -10810
$ws.Range("N65").Value = -19390

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2975.7544
$ws.Range("I132").Value = 3906.5833
$ws.Range("J132").Value = 1380.0476
$ws.Range("K132").Value = 11719.7499
$ws.Range("L132").Value = 4140.142800000001
$ws.Range("M132").Value = -9189.749899999999
$ws.Range("N132").Value = -9200.142800000001
